$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column L (2020) mirroring column K values for each row, copying
# the style used in column K so the new cells look consistent with the
# rest of the table.
$ws.Range("L4").Value = 2020
$ws.Range("L5").Value = 5.6
$ws.Range("L6").Value = 0.8
$ws.Range("L7").Value = 1.9
$ws.Range("L8").Value = 0.7
$ws.Range("L9").Value = 0.7
$ws.Range("L10").Value = 0.9
$ws.Range("L11").Value = 0.3
$ws.Range("L12").Value = 0.2

# Copy formatting from column K into the new column L so styles match.
$ws.Range("K4:K12").Copy() | Out-Null
$ws.Range("L4:L12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Update the active cell selection to match the target workbook.
$ws.Range("N5").Select() | Out-Null
